$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (reordered + new trailing audit columns) ---
$ws.Range("A1").Value = "hierarchy_level"
$ws.Range("B1").Value = "hierarchy_level_name"
$ws.Range("C1").Value = "lang_code"
$ws.Range("D1").Value = "is_active"
$ws.Range("E1").Value = "cr_by"
$ws.Range("F1").Value = "cr_dtimes"
$ws.Range("G1").Value = "upd_by"
$ws.Range("H1").Value = "upd_dtimes"
$ws.Range("I1").Value = "is_deleted"
$ws.Range("J1").Value = "del_dtimes"

# --- Data rows ---
$levels = @("PAYS", "REGION", "PREFECTURE", "SOUS_PREFECTURE_OU_COMMUNE", "DISTRICT", "SECTEUR")
$createdDate = 45079.577674224536

for ($i = 0; $i -lt $levels.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $i
    $ws.Cells.Item($r, 2).Value = $levels[$i]
    $ws.Cells.Item($r, 3).Value = "fra"
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = "superadmin"
    $ws.Cells.Item($r, 6).Value = $createdDate
    $ws.Cells.Item($r, 6).NumberFormat = "mm:ss.0"
    $ws.Cells.Item($r, 7).Value = "NULL"
    $ws.Cells.Item($r, 8).Value = "NULL"
    $ws.Cells.Item($r, 9).Value = $false
    $ws.Cells.Item($r, 10).Value = "NULL"
}

$ws.Range("D16").Select()
